$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (dSF) value updates - repull data / recalculated values
$updates = @{
    2  = -6
    3  = -1
    4  = -3
    6  = 4
    7  = -6
    8  = -9
    13 = -5
    14 = -1
    15 = -2
    16 = 12
    18 = 6
    20 = 0
    23 = -4
    32 = -3
    33 = 0
    34 = -1
    38 = -2
    40 = 6
    42 = -7
    43 = 1
    45 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
